$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "No mapping found for HTTP request with URI [/muziqhub/] in DispatcherServlet with name 'dispatcher'"
$ws.Range("B5").Value = "Project properties-> Facets->Runtime-> add tomcat"

$ws.Range("A6").Value = "Attribute item invalid for tag forEach according to TLD"
$ws.Range("B6").Value = "Add jstl dependencies in pom.xml"

$ws.Range("B8").Select()
